$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Colocação"

$ws.Range("E2").Value = "1º"
$ws.Range("E3").Value = "2º"
$ws.Range("E4").Value = "3º"
$ws.Range("E5").Value = "4º"
$ws.Range("E6").Value = "5º"
$ws.Range("E7").Value = "6º"
$ws.Range("E8").Value = "7º"
